$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they are stored as
# text (matching the source data which is inline-string, not numeric) and
# preserve exact formatting such as trailing zeros (e.g. "73.00").
$ws.Range("D2").Value = '43.355.24'
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = '2.248.95'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '230.46'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '64.10'
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.438'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0952'
$ws.Range("E10").Value = '  -7.97%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '56.84'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '26.62'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '2.584.22'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '14.90'
$ws.Range("E15").Value = '  -5.23%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.819'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '2.249.76'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").Value = '43.239.64'
$ws.Range("E19").Value = '  -1.69%  '
$ws.Range("E20").Value = '  -4.70%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '73.00'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.05'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '246.22'
$ws.Range("E23").Value = '  -3.55%  '
$ws.Range("B24").Value = 'WEMIXToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '3.82'
$ws.Range("E24").Value = '  +14.45%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '2.22'
$ws.Range("E27").Value = '  -4.93%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '9.71'
$ws.Range("E28").Value = '  -4.31%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '173.53'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("E30").Value = '  +3.65%  '
$ws.Range("E31").Value = '  +3.46%  '
$ws.Range("E32").Value = '  -5.57%  '
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '4.90'
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '0.0676'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("E36").Value = '  +1.15%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '3.64'
$ws.Range("E37").Value = '  -5.22%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '6.35'
$ws.Range("E38").Value = '  -5.48%  '
$ws.Range("E39").Value = '  -2.69%  '
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +5.28%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '4.52'
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '17.18'
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("E45").Value = '  +2.62%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '96.52'
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '1.17'
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.0935'
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.000205'
$ws.Range("E49").Value = '  -2.85%  '
$ws.Range("D50").Value = '1.428.89'
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("E51").Value = '  +0.18%  '
